# Rename the "_old" / "_new" header labels to "_FV2404" / "_FV2410"
# and turn the data range A1:U64 into a proper Excel Table (ListObject)
# with a frozen header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename shared-string header cells in row 1 (columns A:U)
$headerRange = $ws.Range("A1:U1")
foreach ($cell in $headerRange.Cells) {
    $text = $cell.Value()
    if ($text -ne $null -and $text -like "*_old") {
        $cell.Value = ($text -replace "_old$", "_FV2404")
    } elseif ($text -ne $null -and $text -like "*_new") {
        $cell.Value = ($text -replace "_new$", "_FV2410")
    }
}

# Convert the data range into an Excel Table named "Table1"
$tableRange = $ws.Range("A1:U64")
$listObject = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$listObject.Name = "Table1"

# Freeze the header row (pane split after row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
